$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the old "_GoBack" bookmark (it used to sit right after
#    " untuk kedepannya"). We delete it now and re-add it later at its
#    new location (end of the "...lebih baik." paragraph), mirroring
#    the diff exactly (a single _GoBack bookmark just moved place).
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2) Split the run containing
#      "Selain itu kafe tersebut juga berharap dapat membuat data
#       pelanggan untuk memberikan pelayanan yang lebih baik."
#    into two runs -- "...lebih bai" and "k." -- where the new,
#    trailing "k." run gets its own rPr (rFonts hint=default,
#    ascii/hAnsi/cs=Tahoma, sz/szCs=24, lang=en-US). Immediately
#    after that new run we (re)place the "_GoBack" bookmark, right
#    before the closing </w:p>.
# ------------------------------------------------------------------
$target = $d.Paragraphs | Where-Object {
    $_.Range.Text -like "*Selain itu kafe tersebut juga berharap dapat membuat data pelanggan untuk memberikan pelayanan yang lebih baik.*"
}

$rng = $target.Range
# Sub-range covering the last two visible characters "k." (End-1 skips
# the paragraph mark itself).
$subRng = $d.Range($rng.End - 3, $rng.End - 1)

$xml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:rFonts w:hint="default" w:ascii="Tahoma" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t>k.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$subRng.InsertXML($xml)

Write-Output "Done"
